$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns are treated as text so numeric-looking
# strings (e.g. "0.999", "90.814.21") are preserved exactly as text
# rather than being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "90.814.21"
$ws.Range("E2").Value = "  +3.80%  "

$ws.Range("D3").Value = "3.203.04"
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "222.33"
$ws.Range("E5").Value = "  +7.52%  "

$ws.Range("D6").Value = "641.49"
$ws.Range("E6").Value = "  +5.54%  "

$ws.Range("D7").Value = "0.401"
$ws.Range("E7").Value = "  +6.28%  "

$ws.Range("D8").Value = "0.709"
$ws.Range("E8").Value = "  +7.00%  "

$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").Value = "3.197.49"
$ws.Range("E10").Value = "  +1.05%  "

$ws.Range("D11").Value = "0.577"
$ws.Range("E11").Value = "  +8.14%  "

$ws.Range("D12").Value = "0.181"
$ws.Range("E12").Value = "  +3.14%  "

$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").Value = "  +7.94%  "

$ws.Range("D14").Value = "5.45"
$ws.Range("E14").Value = "  +4.32%  "

$ws.Range("D15").Value = "33.56"
$ws.Range("E15").Value = "  +4.85%  "

$ws.Range("D16").Value = "90.324.66"
$ws.Range("E16").Value = "  +3.48%  "

$ws.Range("D17").Value = "3.790.95"
$ws.Range("E17").Value = "  +1.08%  "

$ws.Range("D18").Value = "3.206.81"
$ws.Range("E18").Value = "  +0.96%  "

$ws.Range("B19").Value = "PEPE"
$ws.Range("C19").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D19").Value = "0.0000229"
$ws.Range("E19").Value = "  +75.81%  "

$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").Value = "3.36"
$ws.Range("E20").Value = "  +9.06%  "

$ws.Range("D21").Value = "13.50"
$ws.Range("E21").Value = "  +1.08%  "

$ws.Range("D22").Value = "439.69"
$ws.Range("E22").Value = "  +6.60%  "

$ws.Range("D23").Value = "8.65"
$ws.Range("E23").Value = "  +2.81%  "

$ws.Range("D24").Value = "5.08"
$ws.Range("E24").Value = "  +0.96%  "

$ws.Range("D25").Value = "5.38"
$ws.Range("E25").Value = "  +4.62%  "

$ws.Range("D26").Value = "11.92"
$ws.Range("E26").Value = "  +1.15%  "

$ws.Range("D27").Value = "81.50"
$ws.Range("E27").Value = "  +11.77%  "

$ws.Range("D28").Value = "3.373.40"
$ws.Range("E28").Value = "  +1.14%  "

$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").Value = "0.159"
$ws.Range("E30").Value = "  +0.72%  "

$ws.Range("B31").Value = "dogwifhat"
$ws.Range("C31").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D31").Value = "4.29"
$ws.Range("E31").Value = "  +42.82%  "

$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "0.997"
$ws.Range("E32").Value = "  -0.37%  "

$ws.Range("D33").Value = "8.49"
$ws.Range("E33").Value = "  +3.69%  "

$ws.Range("D34").Value = "541.66"
$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("D35").Value = "7.11"
$ws.Range("E35").Value = "  +7.12%  "

$ws.Range("D36").Value = "1.92"
$ws.Range("E36").Value = "  +4.31%  "

$ws.Range("D37").Value = "1.31"
$ws.Range("E37").Value = "  +1.31%  "

$ws.Range("D38").Value = "22.56"
$ws.Range("E38").Value = "  +3.88%  "

$ws.Range("D39").Value = "22.37"
$ws.Range("E39").Value = "  +2.50%  "

$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("D41").Value = "0.127"
$ws.Range("E41").Value = "  -3.82%  "

$ws.Range("D42").Value = "1.96"
$ws.Range("E42").Value = "  +2.95%  "

$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("D44").Value = "0.375"
$ws.Range("E44").Value = "  +2.30%  "

$ws.Range("D45").Value = "146.76"
$ws.Range("E45").Value = "  -1.46%  "

$ws.Range("D46").Value = "44.89"
$ws.Range("E46").Value = "  +4.25%  "

$ws.Range("D47").Value = "173.50"
$ws.Range("E47").Value = "  +0.82%  "

$ws.Range("D48").Value = "0.126"
$ws.Range("E48").Value = "  +1.84%  "

$ws.Range("D49").Value = "0.751"
$ws.Range("E49").Value = "  +8.70%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "25.09"
$ws.Range("E50").Value = "  +6.85%  "

$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "0.624"
$ws.Range("E51").Value = "  +7.28%  "

# Restore default (General) number format and default style so the
# untouched formatting metadata matches the original workbook as
# closely as possible, while keeping the values stored as text.
$ws.Range("D2:E51").NumberFormat = "General"
$ws.Range("D2:E51").Style = "Normal"
